$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.958.07"
$ws.Range("E2").Value = "  -4.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.272.38"
$ws.Range("E3").Value = "  -5.73%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.94"
$ws.Range("E5").Value = "  -3.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.40"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.262.88"
$ws.Range("E9").Value = "  -5.76%  "
$ws.Range("E10").Value = "  -8.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.587"
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.37"
$ws.Range("E13").Value = "  -6.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.59"
$ws.Range("E14").Value = "  -5.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "629.77"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.802.77"
$ws.Range("E16").Value = "  -5.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.877.88"
$ws.Range("E17").Value = "  -4.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.86"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.280.27"
$ws.Range("E20").Value = "  -5.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.33"
$ws.Range("E21").Value = "  -7.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.903"
$ws.Range("E22").Value = "  -4.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.15"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "106.42"
$ws.Range("E24").Value = "  +7.45%  "
$ws.Range("E25").Value = "  -7.04%  "
$ws.Range("E26").Value = "  -7.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("E27").Value = "  -6.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.70"
$ws.Range("E29").Value = "  -6.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.34"
$ws.Range("E30").Value = "  -6.13%  "
$ws.Range("E31").Value = "  -6.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  -6.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.03"
$ws.Range("E33").Value = "  -4.78%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.105"
$ws.Range("E34").Value = "  -3.76%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "537.27"
$ws.Range("E35").Value = "  +5.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.45"
$ws.Range("E36").Value = "  -5.75%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.718.78"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0728"
$ws.Range("E40").Value = "  -7.95%  "
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("E42").Value = "  -6.81%  "
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.78"
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.338"
$ws.Range("E45").Value = "  -8.90%  "
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.60"
$ws.Range("E49").Value = "  -7.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  +2.16%  "
